# Update countries & provincias Spain
# Applies the data refresh captured by the diff between the previous and
# updated "paises.xlsx" COVID tracking workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 02:20"

# 2. Estados Unidos (row 4): refresh case counters
$ws.Cells.Item(4, 2).Value = 141812   # Casos totales
$ws.Cells.Item(4, 3).Value = 18234    # Nuevos casos
$ws.Cells.Item(4, 5).Value = 134902   # Recuperados
$ws.Cells.Item(4, 7).Value = 255      # Muertes hoy
$ws.Cells.Item(4, 8).Value = 2475     # Muertes

# 3. Iran (row 8): refresh case counters
$ws.Cells.Item(8, 5).Value = 52343    # Recuperados
$ws.Cells.Item(8, 7).Value = 108      # Muertes hoy
$ws.Cells.Item(8, 8).Value = 541      # Muertes

# 4. Vietnam (row 91): refresh case counters
$ws.Cells.Item(91, 2).Value = 194     # Casos totales
$ws.Cells.Item(91, 3).Value = 20      # Nuevos casos
$ws.Cells.Item(91, 5).Value = 169     # Recuperados

# 5. Reorder "Consejo Danes para los Refugiados" so it appears right after
#    Bolivia (row 117) instead of after Ruanda, and refresh its counters.
#    Insert a fresh row at 118 (this pushes Trinidad yTobago and Ruanda
#    down by one row, to 119 and 120) and populate it with the updated
#    figures for the Consejo Danes entry.
$ws.Rows(118).Insert()

$ws.Cells.Item(118, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(118, 2).Value = 81
$ws.Cells.Item(118, 3).Value = 16
$ws.Cells.Item(118, 4).Value = 2
$ws.Cells.Item(118, 5).Value = 71
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 2
$ws.Cells.Item(118, 8).Value = 8

# Now delete the old "Consejo Danes para los Refugiados" row, which has
# been shifted down to row 121 (Trinidad yTobago=119, Ruanda=120,
# old Consejo Danes=121, Gibraltar=122).
$ws.Rows(121).Delete()
